$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to be treated as plain text so the dashed dates are
# stored as literal strings (matching the inlineStr cells in the diff)
# rather than being auto-parsed into Excel date serials.
$dateCells = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $dateCells) {
    $ws.Range("A$r").NumberFormat = "@"
}

# Update dates (slash -> dash) and D/E/G/H counts for rows 3-21
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("A4").Value = "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("A5").Value = "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("A6").Value = "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

$ws.Range("A7").Value = "11-08-2022"

$ws.Range("A8").Value = "15-08-2022"

$ws.Range("A9").Value = "18-08-2022"

$ws.Range("A10").Value = "22-08-2022"

$ws.Range("A11").Value = "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

$ws.Range("A12").Value = "29-08-2022"

$ws.Range("A13").Value = "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

$ws.Range("A14").Value = "05-09-2022"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

$ws.Range("A15").Value = "08-09-2022"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("H15").Value = 0

$ws.Range("A16").Value = "12-09-2022"

$ws.Range("A17").Value = "15-09-2022"

$ws.Range("A18").Value = "19-09-2022"

$ws.Range("A19").Value = "22-09-2022"

$ws.Range("A20").Value = "26-09-2022"

$ws.Range("A21").Value = "29-09-2022"
